$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.731.18"
$ws.Range("E2").Value = "  -3.82%  "

$ws.Range("D3").Value = "'2.981.04"
$ws.Range("E3").Value = "  -4.86%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'545.03"
$ws.Range("E5").Value = "  -4.28%  "

$ws.Range("D6").Value = "'153.22"
$ws.Range("E6").Value = "  -5.02%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("D9").Value = "'2.991.01"
$ws.Range("E9").Value = "  -4.86%  "

$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -6.16%  "

$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "'3.503.76"
$ws.Range("E13").Value = "  -4.84%  "

$ws.Range("D14").Value = "'0.125"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").Value = "'61.797.21"
$ws.Range("E15").Value = "  -3.84%  "

$ws.Range("E16").Value = "  -4.18%  "

$ws.Range("D17").Value = "'2.986.05"
$ws.Range("E17").Value = "  -4.97%  "

$ws.Range("D18").Value = "'0.0000148"
$ws.Range("E18").Value = "  -3.48%  "

$ws.Range("D19").Value = "'5.17"
$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'383.25"
$ws.Range("E20").Value = "  -4.36%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'12.04"
$ws.Range("E21").Value = "  -3.42%  "

$ws.Range("E22").Value = "  -5.24%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'65.91"
$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("E25").Value = "  -1.61%  "

$ws.Range("D26").Value = "'3.097.07"
$ws.Range("E26").Value = "  -5.47%  "

$ws.Range("D27").Value = "'0.189"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("D28").Value = "'0.0₃0948"
$ws.Range("E28").Value = "  -5.31%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").Value = "'8.31"
$ws.Range("E30").Value = "  -5.20%  "

$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("E32").Value = "  -3.69%  "

$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").Value = "'160.87"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").Value = "'4.70"
$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").Value = "'5.99"
$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("E37").Value = "  -1.71%  "

$ws.Range("D38").Value = "'1.28"
$ws.Range("E38").Value = "  -4.03%  "

$ws.Range("E39").Value = "  -5.44%  "

$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("D41").Value = "'2.414.16"
$ws.Range("E41").Value = "  -9.19%  "

$ws.Range("D42").Value = "'37.42"
$ws.Range("E42").Value = "  -2.29%  "

$ws.Range("D43").Value = "'22.33"
$ws.Range("E43").Value = "  -4.46%  "

$ws.Range("D44").Value = "'0.669"
$ws.Range("E44").Value = "  -2.37%  "

$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0249"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'5.10"
$ws.Range("E47").Value = "  -5.76%  "

$ws.Range("D48").Value = "'0.996"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'272.13"
$ws.Range("E49").Value = "  -5.14%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'19.96"
$ws.Range("E50").Value = "  -4.57%  "

$ws.Range("D51").Value = "'0.0953"
$ws.Range("E51").Value = "  -1.70%  "

